$wb = $excel.ActiveWorkbook

# --- Fix typo in the BDCTBA sheet label (B1): "Diable" -> "Disable" ---
$wsBDCTBA = $wb.Worksheets.Item("BDCTBA")
$wsBDCTBA.Range("B1").Value = "Disable Carbon Tax Border Adjustment"

# --- Set the actual boolean lever value to 1 (disable CTBA) ---
$wsBDCTBA.Range("B2").Value = 1

# --- Selection / active-sheet bookkeeping to match the authored state ---
$wsBDCTBA.Range("B2").Select()

$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Activate()
$wsAbout.Range("B15").Select()
